$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.20014
$ws.Range("H2").Value = 0.6004200000000001
$ws.Range("I2").Value = 0.09409628186320101
$ws.Range("J2").Value = 0.09409628186320104
$ws.Range("M2").Value = 73.202511
$ws.Range("N2").Value = 219.607533
$ws.Range("O2").Value = 0.3264904632507938
$ws.Range("P2").Value = 0.3264904632507938
$ws.Range("Q2").Value = 14.65075055154
$ws.Range("R2").Value = 131.85675496386
$ws.Range("S2").Value = 0.03072153865569376
$ws.Range("T2").Value = 0.03072153865569378
$ws.Range("G3").Value = 0.20014
$ws.Range("H3").Value = 0.6004200000000001
$ws.Range("I3").Value = 0.09409628186320101
$ws.Range("J3").Value = 0.09409628186320104
$ws.Range("O3").Value = 0.4449719839907295
$ws.Range("P3").Value = 0.4449719839907295
$ws.Range("Q3").Value = 19.96742408633333
$ws.Range("R3").Value = 179.706816777
$ws.Range("S3").Value = 0.04187020922681945
$ws.Range("T3").Value = 0.04187020922681946
$ws.Range("G4").Value = 0.20014
$ws.Range("H4").Value = 0.6004200000000001
$ws.Range("I4").Value = 0.09409628186320101
$ws.Range("J4").Value = 0.09409628186320104
$ws.Range("M4").Value = 39.54025133333334
$ws.Range("N4").Value = 118.620754
$ws.Range("O4").Value = 0.1763534446908907
$ws.Range("P4").Value = 0.1763534446908907
$ws.Range("Q4").Value = 7.913585901853335
$ws.Range("R4").Value = 71.22227311668001
$ws.Range("S4").Value = 0.01659420343918048
$ws.Range("T4").Value = 0.01659420343918049
$ws.Range("G5").Value = 0.20014
$ws.Range("H5").Value = 0.6004200000000001
$ws.Range("I5").Value = 0.09409628186320101
$ws.Range("J5").Value = 0.09409628186320104
$ws.Range("M5").Value = 11.70021233333333
$ws.Range("N5").Value = 35.100637
$ws.Range("O5").Value = 0.05218410806758597
$ws.Range("P5").Value = 0.05218410806758598
$ws.Range("Q5").Value = 2.341680496393333
$ws.Range("R5").Value = 21.07512446754
$ws.Range("S5").Value = 0.004910330541507311
$ws.Range("T5").Value = 0.004910330541507314
$ws.Range("I6").Value = 0.7283659026117116
$ws.Range("J6").Value = 0.7283659026117117
$ws.Range("M6").Value = 73.202511
$ws.Range("N6").Value = 219.607533
$ws.Range("O6").Value = 0.3264904632507938
$ws.Range("P6").Value = 0.3264904632507938
$ws.Range("Q6").Value = 113.406257273006
$ws.Range("R6").Value = 1020.656315457054
$ws.Range("S6").Value = 0.2378045209597803
$ws.Range("T6").Value = 0.2378045209597804
$ws.Range("I7").Value = 0.7283659026117116
$ws.Range("J7").Value = 0.7283659026117117
$ws.Range("O7").Value = 0.4449719839907295
$ws.Range("P7").Value = 0.4449719839907295
$ws.Range("S7").Value = 0.3241024207563318
$ws.Range("T7").Value = 0.3241024207563318
$ws.Range("I8").Value = 0.7283659026117116
$ws.Range("J8").Value = 0.7283659026117117
$ws.Range("M8").Value = 39.54025133333334
$ws.Range("N8").Value = 118.620754
$ws.Range("O8").Value = 0.1763534446908907
$ws.Range("P8").Value = 0.1763534446908907
$ws.Range("Q8").Value = 61.25625820878356
$ws.Range("R8").Value = 551.3063238790521
$ws.Range("S8").Value = 0.1284498359209652
$ws.Range("T8").Value = 0.1284498359209652
$ws.Range("I9").Value = 0.7283659026117116
$ws.Range("J9").Value = 0.7283659026117117
$ws.Range("M9").Value = 11.70021233333333
$ws.Range("N9").Value = 35.100637
$ws.Range("O9").Value = 0.05218410806758597
$ws.Range("P9").Value = 0.05218410806758598
$ws.Range("Q9").Value = 18.12611714948956
$ws.Range("R9").Value = 163.135054345406
$ws.Range("S9").Value = 0.03800912497463436
$ws.Range("T9").Value = 0.03800912497463437
$ws.Range("G10").Value = 0.25539
$ws.Range("H10").Value = 0.76617
$ws.Range("I10").Value = 0.1200721965876032
$ws.Range("J10").Value = 0.1200721965876032
$ws.Range("M10").Value = 73.202511
$ws.Range("N10").Value = 219.607533
$ws.Range("O10").Value = 0.3264904632507938
$ws.Range("P10").Value = 0.3264904632507938
$ws.Range("Q10").Value = 18.69518928429
$ws.Range("R10").Value = 168.25670355861
$ws.Range("S10").Value = 0.03920242708742696
$ws.Range("T10").Value = 0.03920242708742697
$ws.Range("G11").Value = 0.25539
$ws.Range("H11").Value = 0.76617
$ws.Range("I11").Value = 0.1200721965876032
$ws.Range("J11").Value = 0.1200721965876032
$ws.Range("O11").Value = 0.4449719839907295
$ws.Range("P11").Value = 0.4449719839907295
$ws.Range("Q11").Value = 25.4795664905
$ws.Range("R11").Value = 229.3160984145
$ws.Range("S11").Value = 0.05342876353771071
$ws.Range("T11").Value = 0.05342876353771071
$ws.Range("G12").Value = 0.25539
$ws.Range("H12").Value = 0.76617
$ws.Range("I12").Value = 0.1200721965876032
$ws.Range("J12").Value = 0.1200721965876032
$ws.Range("M12").Value = 39.54025133333334
$ws.Range("N12").Value = 118.620754
$ws.Range("O12").Value = 0.1763534446908907
$ws.Range("P12").Value = 0.1763534446908907
$ws.Range("Q12").Value = 10.09818478802
$ws.Range("R12").Value = 90.88366309218
$ws.Range("S12").Value = 0.02117514547982564
$ws.Range("T12").Value = 0.02117514547982564
$ws.Range("G13").Value = 0.25539
$ws.Range("H13").Value = 0.76617
$ws.Range("I13").Value = 0.1200721965876032
$ws.Range("J13").Value = 0.1200721965876032
$ws.Range("M13").Value = 11.70021233333333
$ws.Range("N13").Value = 35.100637
$ws.Range("O13").Value = 0.05218410806758597
$ws.Range("P13").Value = 0.05218410806758598
$ws.Range("Q13").Value = 2.98811722781
$ws.Range("R13").Value = 26.89305505029
$ws.Range("S13").Value = 0.006265860482639914
$ws.Range("T13").Value = 0.006265860482639916
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1222276666666667
$ws.Range("H14").Value = 0.366683
$ws.Range("I14").Value = 0.05746561893748399
$ws.Range("J14").Value = 0.057465618937484
$ws.Range("M14").Value = 73.202511
$ws.Range("N14").Value = 219.607533
$ws.Range("O14").Value = 0.3264904632507938
$ws.Range("P14").Value = 0.3264904632507938
$ws.Range("Q14").Value = 8.947372113670999
$ws.Range("R14").Value = 80.52634902303899
$ws.Range("S14").Value = 0.01876197654789273
$ws.Range("T14").Value = 0.01876197654789274
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1222276666666667
$ws.Range("H15").Value = 0.366683
$ws.Range("I15").Value = 0.05746561893748399
$ws.Range("J15").Value = 0.057465618937484
$ws.Range("O15").Value = 0.4449719839907295
$ws.Range("P15").Value = 0.4449719839907295
$ws.Range("Q15").Value = 12.19432225150556
$ws.Range("R15").Value = 109.74890026355
$ws.Range("S15").Value = 0.02557059046986749
$ws.Range("T15").Value = 0.0255705904698675
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1222276666666667
$ws.Range("H16").Value = 0.366683
$ws.Range("I16").Value = 0.05746561893748399
$ws.Range("J16").Value = 0.057465618937484
$ws.Range("M16").Value = 39.54025133333334
$ws.Range("N16").Value = 118.620754
$ws.Range("O16").Value = 0.1763534446908907
$ws.Range("P16").Value = 0.1763534446908907
$ws.Range("Q16").Value = 4.832912659886889
$ws.Range("R16").Value = 43.496213938982
$ws.Range("S16").Value = 0.01013425985091938
$ws.Range("T16").Value = 0.01013425985091939
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1222276666666667
$ws.Range("H17").Value = 0.366683
$ws.Range("I17").Value = 0.05746561893748399
$ws.Range("J17").Value = 0.057465618937484
$ws.Range("M17").Value = 11.70021233333333
$ws.Range("N17").Value = 35.100637
$ws.Range("O17").Value = 0.05218410806758597
$ws.Range("P17").Value = 0.05218410806758598
$ws.Range("Q17").Value = 1.430089653007889
$ws.Range("R17").Value = 12.870806877071
$ws.Range("S17").Value = 0.002998792068804379
$ws.Range("T17").Value = 0.002998792068804381
